# Updates the cryptocurrency "Price" (D) and "Volume(1h)" (E) columns on
# the active worksheet with the latest scrape, per the GitHub Actions
# "Updated cryptos list" commit. Row 23 (Dai) is unchanged in the source
# data, so it is left untouched here.
#
# Price cells are forced to Text format ("@") before the assignment so
# numeric-looking strings (e.g. "536.48") round-trip as text instead of
# being auto-coerced to a Number by the COM Value setter - matching the
# original workbook, where these cells are stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.651.96"
$ws.Range("E2").Value = "  -1.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.633.18"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.48"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.20"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.642.27"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.04"
$ws.Range("E10").Value = "  +8.46%  "
$ws.Range("E11").Value = "  -1.64%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.096.01"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.597.02"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.91"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.639.61"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.16"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.17"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.38"
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.164"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.16"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0740"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.84"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.76"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.46"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.91"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "37.13"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.11"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.847"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("E39").Value = "  -2.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.815"
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "282.43"
$ws.Range("E42").Value = "  +3.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.70"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.15"
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0936"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.949.03"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.51"
$ws.Range("E51").Value = "  +0.03%  "
